$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 340 has the exact same column layout as the rows we are appending
# (A type / B date(style1) / C "Global" / D "M" centered(style3) / E player /
#  F poste / G temps joue / H..V numeric stats), so use it as a formatting
# template for the new rows.
$templateRow = 340

$players = @(
    @{Row=1038; Poste="right back";      Temps="01:38:41"; H=10.82; I=1.45; J=9.35;               K=0.95; L=0.32;               M=0.18; N=0.01; O=11; P=6.43; Q=31.05; R=4.79;               S=39; T=9;  U=35; V=17},
    @{Row=1039; Poste="left back";       Temps="01:39:40"; H=11.78; I=1.64; J=10.11;              K=1.25; L=0.38;               M=0.03; N=0;    O=4;  P=6.87; Q=26.59; R=4.79;               S=46; T=12; U=46; V=13},
    @{Row=1040; Poste="left forward";    Temps="01:27:54"; H=10.84; I=2.44; J=8.3699999999999992; K=1.46; L=0.74;               M=0.24; N=0.02; O=16; P=7.21; Q=30.89; R=4.91;               S=58; T=11; U=39; V=20},
    @{Row=1041; Poste="right forward";   Temps="01:33:05"; H=9.65;  I=1.38; J=8.25;               K=1;    L=0.22;               M=0.09; N=0.08; O=7;  P=6.04; Q=32.33; R=4.97;               S=40; T=9;  U=27; V=13},
    @{Row=1042; Poste="center back";     Temps="01:39:23"; H=10.28; I=1.35; J=8.91;               K=0.92; L=0.28000000000000003; M=0.13; N=0.05; O=11; P=6.07; Q=31.44; R=4.47;               S=41; T=5;  U=30; V=8},
    @{Row=1043; Poste="center midfield"; Temps="01:39:23"; H=11.25; I=1.96; J=9.27;               K=1.42; L=0.48;               M=0.06; N=0.03; O=5;  P=6.76; Q=31.59; R=4.8499999999999996; S=33; T=5;  U=37; V=9},
    @{Row=1044; Poste="left forward";    Temps="01:25:14"; H=9.6999999999999993; I=1.54; J=8.1300000000000008; K=1.07; L=0.39; M=0.11; N=0;    O=8;  P=6.81; Q=29.85; R=4.57;               S=26; T=4;  U=16; V=17},
    @{Row=1045; Poste="center back";     Temps="01:38:41"; H=10.64; I=1.63; J=8.99;               K=1.24; L=0.3;                M=0.1;  N=0;    O=7;  P=6.34; Q=28.9;  R=4.5;                S=35; T=4;  U=42; V=5},
    @{Row=1046; Poste="center back";     Temps="01:26:21"; H=8.6999999999999993; I=0.83; J=7.86; K=0.55000000000000004; L=0.16; M=0.12; N=0.01; O=7;  P=5.75; Q=31.36; R=4.34;               S=22; T=1;  U=17; V=5},
    @{Row=1047; Poste="center midfield"; Temps="01:37:42"; H=12.13; I=2.64; J=9.4499999999999993; K=2.02; L=0.55000000000000004; M=0.11; N=0; O=6;  P=7.42; Q=28.28; R=4.24;               S=19; T=1;  U=38; V=6},
    @{Row=1048; Poste="center midfield"; Temps="00:13:53"; H=1.82;  I=0.53; J=1.29;               K=0.39; L=0.12;               M=0.02; N=0;    O=1;  P=7.82; Q=30.31; R=4.08;               S=10; T=1;  U=1;  V=1}
)
$names = @("Mattheo Haon", "Maé Clavel", "Emmanuel Valey", "Amir Etien", "Naim Ighbane", "Naim Dhib", "Sofiane Belle", "Romain Thunet", "Yoan Zouma", "Yoann Martelat", "Ilan Ihaddadene")

for ($idx = 0; $idx -lt $players.Count; $idx++) {
    $p = $players[$idx]
    $r = $p.Row

    # Bring over formatting (column widths/styles) from the template row.
    $ws.Range("A$templateRow`:V$templateRow").Copy()
    $ws.Range("A$r`:V$r").PasteSpecial(-4122)

    $ws.Cells.Item($r, 2).Value = 46011
    $ws.Cells.Item($r, 3).Value = "Global"
    $ws.Cells.Item($r, 4).Value = "M"
    $ws.Cells.Item($r, 5).Value = $names[$idx]
    $ws.Cells.Item($r, 6).Value = $p.Poste
}

# Column G ("Temps joué") filled top-to-bottom across the new rows, mirroring
# the order the shared-string table records these newly-seen time values.
for ($idx = 0; $idx -lt $players.Count; $idx++) {
    $p = $players[$idx]
    $ws.Cells.Item($p.Row, 7).Value = $p.Temps
}

for ($idx = 0; $idx -lt $players.Count; $idx++) {
    $p = $players[$idx]
    $r = $p.Row
    $ws.Cells.Item($r, 8).Value = $p.H
    $ws.Cells.Item($r, 9).Value = $p.I
    $ws.Cells.Item($r, 10).Value = $p.J
    $ws.Cells.Item($r, 11).Value = $p.K
    $ws.Cells.Item($r, 12).Value = $p.L
    $ws.Cells.Item($r, 13).Value = $p.M
    $ws.Cells.Item($r, 14).Value = $p.N
    $ws.Cells.Item($r, 15).Value = $p.O
    $ws.Cells.Item($r, 16).Value = $p.P
    $ws.Cells.Item($r, 17).Value = $p.Q
    $ws.Cells.Item($r, 18).Value = $p.R
    $ws.Cells.Item($r, 19).Value = $p.S
    $ws.Cells.Item($r, 20).Value = $p.T
    $ws.Cells.Item($r, 21).Value = $p.U
    $ws.Cells.Item($r, 22).Value = $p.V
}

# Column A (Type / match description) is identical for every new row, and is
# the last new shared string created - fill it in after every other column.
for ($idx = 0; $idx -lt $players.Count; $idx++) {
    $p = $players[$idx]
    $ws.Cells.Item($p.Row, 1).Value = "CDF 32ème VS Toulouse FC (L1) "
}

$ws.Application.ActiveWindow.ScrollRow = 1017
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("E1052").Select()
